$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = -0.6231890030091449
$ws.Range("C2").Value = -0.3815599582472184
$ws.Range("D2").Value = 1.039208715109508
$ws.Range("E2").Value = 1.562018419371699
$ws.Range("B3").Value = 0.1313269621646771
$ws.Range("C3").Value = 0.8263615303988285
$ws.Range("D3").Value = 3.145821501937435
$ws.Range("E3").Value = 3.668631206199626
$ws.Range("B4").Value = 0.343234630204168
$ws.Range("C4").Value = 0.922789407056382
$ws.Range("D4").Value = 4.810881543137723
$ws.Range("E4").Value = 5.333691247399914
$ws.Range("B5").Value = 0.7125522396455506
$ws.Range("C5").Value = 0.4815132279214041
$ws.Range("D5").Value = 2.567093677273756
$ws.Range("E5").Value = 3.089903381535947
$ws.Range("B6").Value = 0.8662106395267892
$ws.Range("C6").Value = 0.1848282929284375
$ws.Range("D6").Value = 2.551152949613654
$ws.Range("E6").Value = 3.073962653875845
$ws.Range("B7").Value = -0.7743300267179998
$ws.Range("C7").Value = 0.07101619575450724
$ws.Range("D7").Value = 1.675734724906242
$ws.Range("E7").Value = 2.198544429168433
$ws.Range("B8").Value = 0.6048556880868585
$ws.Range("C8").Value = 0.6871356504909893
$ws.Range("D8").Value = 3.311757838029258
$ws.Range("E8").Value = 3.834567542291449
$ws.Range("B9").Value = 0.6311408381149812
$ws.Range("C9").Value = 0.8219340074643868
$ws.Range("D9").Value = 4.7030298689657
$ws.Range("E9").Value = 5.225839573227891
$ws.Range("B10").Value = -0.2162434898902577
$ws.Range("C10").Value = -0.8024443594015642
$ws.Range("D10").Value = -2.307207602449695
$ws.Range("E10").Value = -1.784397898187504
$ws.Range("B11").Value = 0.2977985744636613
$ws.Range("C11").Value = -0.4769871861722883
$ws.Range("D11").Value = -0.4137742482290986
$ws.Range("E11").Value = 0.1090354560330925
$ws.Range("B12").Value = 0.623077212741084
$ws.Range("C12").Value = 0.4963231194350934
$ws.Range("D12").Value = 2.221809807309105
$ws.Range("E12").Value = 2.744619511571297
$ws.Range("B13").Value = 0.5205295124310321
$ws.Range("C13").Value = 0.09129625313229872
$ws.Range("D13").Value = 0.926196165218044
$ws.Range("E13").Value = 1.449005869480235
$ws.Range("B14").Value = -0.916138814832975
$ws.Range("C14").Value = 0.2995899633392989
$ws.Range("D14").Value = 2.325232726393821
$ws.Range("E14").Value = 2.848042430656012
$ws.Range("B15").Value = -0.6876339652659114
$ws.Range("C15").Value = -0.3308805041004506
$ws.Range("D15").Value = 1.374015682954193
$ws.Range("E15").Value = 1.896825387216384
$ws.Range("B16").Value = -0.2690791256533465
$ws.Range("C16").Value = 0.4919033169845433
$ws.Range("D16").Value = 0.7264555764250926
$ws.Range("E16").Value = 1.249265280687284
$ws.Range("B17").Value = -0.4749217012840903
$ws.Range("C17").Value = 0.994134128454679
$ws.Range("D17").Value = 5.267794549091613
$ws.Range("E17").Value = 5.790604253353805
$ws.Range("B18").Value = 0.5228567667030581
$ws.Range("C18").Value = 0.03072945381546499
$ws.Range("D18").Value = 0.8918028393245876
$ws.Range("E18").Value = 1.414612543586779
$ws.Range("B19").Value = -0.05567193566994177
$ws.Range("C19").Value = -0.8750665542320568
$ws.Range("D19").Value = -3.308858726275794
$ws.Range("E19").Value = -2.786049022013603
$ws.Range("B20").Value = -0.5499603893133345
$ws.Range("C20").Value = -0.971532119186785
$ws.Range("D20").Value = -3.201111049687445
$ws.Range("E20").Value = -2.678301345425254
$ws.Range("B21").Value = -0.3722712987954699
$ws.Range("C21").Value = -0.687399935006419
$ws.Range("D21").Value = -1.011104747457896
$ws.Range("E21").Value = -0.4882950431957049
$ws.Range("B22").Value = -0.8036632759173747
$ws.Range("C22").Value = 0.6512951215389555
$ws.Range("D22").Value = 2.822732927983922
$ws.Range("E22").Value = 3.345542632246112
$ws.Range("B23").Value = 0.5773802583911078
$ws.Range("C23").Value = 0.3034691942391716
$ws.Range("D23").Value = 1.412353090447157
$ws.Range("E23").Value = 1.935162794709348
$ws.Range("B24").Value = 0.9406693238294122
$ws.Range("C24").Value = -0.006699004185547475
$ws.Range("D24").Value = 2.741674799129028
$ws.Range("E24").Value = 3.264484503391219
$ws.Range("B25").Value = -0.8146686208556893
$ws.Range("C25").Value = 0.5631427997865806
$ws.Range("D25").Value = 2.431787370368446
$ws.Range("E25").Value = 2.954597074630637
